$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2001.3
$ws.Range("I98").Value = 2651.5
$ws.Range("J98").Value = 1351.1
$ws.Range("K98").Value = 2651.5
$ws.Range("L98").Value = 1351.1
$ws.Range("M98").Value = -1153.5
$ws.Range("N98").Value = -4347.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7181
$ws.Range("I116").Value = 1385
$ws.Range("J116").Value = 15875
$ws.Range("K116").Value = 1385
$ws.Range("L116").Value = 15875
$ws.Range("M116").Value = 2057
$ws.Range("N116").Value = -22759

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 737.5
$ws.Range("I121").Value = 397.5
$ws.Range("J121").Value = 768.4091
$ws.Range("K121").Value = 1192.5
$ws.Range("L121").Value = 2305.2273
$ws.Range("M121").Value = 554.5
$ws.Range("N121").Value = -5799.2273

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2001.3
$ws.Range("I122").Value = 2651.5
$ws.Range("J122").Value = 1351.1
$ws.Range("K122").Value = 7954.5
$ws.Range("L122").Value = 4053.3
$ws.Range("M122").Value = -5504.5
$ws.Range("N122").Value = -8953.299999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2566301.5
$ws.Range("I137").Value = 5557292
$ws.Range("J137").Value = 2595.3809
$ws.Range("K137").Value = 16671876
$ws.Range("L137").Value = 7786.1427
$ws.Range("M137").Value = -16669326
$ws.Range("N137").Value = -12886.1427

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1834912.5
$ws.Range("J138").Value = 2529352.5
$ws.Range("L138").Value = 7588057.5
$ws.Range("N138").Value = -7598337.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 64780
$ws.Range("J140").Value = 64780
$ws.Range("L140").Value = 64780
$ws.Range("N140").Value = -75140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8690763
$ws.Range("I74").Value = 12550796
$ws.Range("J74").Value = 112910.664
$ws.Range("K74").Value = 12550796
$ws.Range("L74").Value = 112910.664
$ws.Range("M74").Value = -12549922
$ws.Range("N74").Value = -114658.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 8690763
$ws.Range("I77").Value = 12550796
$ws.Range("J77").Value = 112910.664
$ws.Range("K77").Value = 62753980
$ws.Range("L77").Value = 564553.3200000001
$ws.Range("M77").Value = -62749612
$ws.Range("N77").Value = -573289.3200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 35000
$ws.Range("J141").Value = 35000
$ws.Range("L141").Value = 35000
$ws.Range("N141").Value = -45360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7307.8477
$ws.Range("I31").Value = 27488.191
$ws.Range("J31").Value = 1339.014
$ws.Range("K31").Value = 27488.191
$ws.Range("L31").Value = 1339.014
$ws.Range("M31").Value = -27193.191
$ws.Range("N31").Value = -1929.014

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7307.8477
$ws.Range("I34").Value = 27488.191
$ws.Range("J34").Value = 1339.014
$ws.Range("K34").Value = 27488.191
$ws.Range("L34").Value = 1339.014
$ws.Range("M34").Value = -27286.191
$ws.Range("N34").Value = -1743.014

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4501.7144
$ws.Range("I99").Value = 4984.1665
$ws.Range("J99").Value = 1607
$ws.Range("K99").Value = 4984.1665
$ws.Range("L99").Value = 1607
$ws.Range("M99").Value = -3486.1665
$ws.Range("N99").Value = -4603

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4501.7144
$ws.Range("I126").Value = 4984.1665
$ws.Range("J126").Value = 1607
$ws.Range("K126").Value = 14952.4995
$ws.Range("L126").Value = 4821
$ws.Range("M126").Value = -12482.4995
$ws.Range("N126").Value = -9761

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 42763.04
$ws.Range("I132").Value = 2541.4546
$ws.Range("J132").Value = 74365.71000000001
$ws.Range("K132").Value = 7624.3638
$ws.Range("L132").Value = 223097.13
$ws.Range("M132").Value = -5094.3638
$ws.Range("N132").Value = -228157.13

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 36272.97
$ws.Range("I134").Value = 1277.4
$ws.Range("J134").Value = 99901.27
$ws.Range("K134").Value = 3832.2
$ws.Range("L134").Value = 299703.81
$ws.Range("M134").Value = -1297.2
$ws.Range("N134").Value = -304773.81

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 787.7586
$ws.Range("J131").Value = 908.8570999999999
$ws.Range("L131").Value = 2726.5713
$ws.Range("N131").Value = -12806.5713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 35142.332
$ws.Range("J137").Value = 48302.46
$ws.Range("L137").Value = 144907.38
$ws.Range("N137").Value = -155107.38

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2209.7693
$ws.Range("I102").Value = 2314.25
$ws.Range("J102").Value = 2042.6
$ws.Range("K102").Value = 2314.25
$ws.Range("L102").Value = 2042.6
$ws.Range("M102").Value = -692.25
$ws.Range("N102").Value = -5286.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2892.75
$ws.Range("I126").Value = 1716.6666
$ws.Range("J126").Value = 4068.8333
$ws.Range("K126").Value = 5149.9998
$ws.Range("L126").Value = 12206.4999
$ws.Range("M126").Value = -2679.9998
$ws.Range("N126").Value = -17146.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2800
$ws.Range("I7").Value = 2700
$ws.Range("J7").Value = 2850
$ws.Range("K7").Value = 2700
$ws.Range("L7").Value = 2850
$ws.Range("M7").Value = -2588
$ws.Range("N7").Value = -3074

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6029.2666
$ws.Range("I40").Value = 6002.8335
$ws.Range("J40").Value = 6135
$ws.Range("K40").Value = 6002.8335
$ws.Range("L40").Value = 6135
$ws.Range("M40").Value = -5866.8335
$ws.Range("N40").Value = -6407

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 600
$ws.Range("I46").Value = 600
$ws.Range("K46").Value = 600
$ws.Range("M46").Value = -412

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2933.3333
$ws.Range("I82").Value = 3000
$ws.Range("J82").Value = 2866.6667
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 2866.6667
$ws.Range("M82").Value = -2639
$ws.Range("N82").Value = -3588.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2933.3333
$ws.Range("I85").Value = 3000
$ws.Range("J85").Value = 2866.6667
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 2866.6667
$ws.Range("M85").Value = -1752
$ws.Range("N85").Value = -5362.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2800
$ws.Range("I126").Value = 2700
$ws.Range("J126").Value = 2850
$ws.Range("K126").Value = 8100
$ws.Range("L126").Value = 8550
$ws.Range("M126").Value = -5630
$ws.Range("N126").Value = -13490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4358.7856
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 4463.3076
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 4463.3076
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -5711.3076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4358.7856
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 4463.3076
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 22316.538
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -28556.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2238.5715
$ws.Range("I126").Value = 1417.5
$ws.Range("J126").Value = 3333.3333
$ws.Range("K126").Value = 4252.5
$ws.Range("L126").Value = 9999.999899999999
$ws.Range("M126").Value = -1782.5
$ws.Range("N126").Value = -14939.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 60870.266
$ws.Range("I132").Value = 41891.48
$ws.Range("J132").Value = 113589.11
$ws.Range("K132").Value = 125674.44
$ws.Range("L132").Value = 340767.33
$ws.Range("M132").Value = -123144.44
$ws.Range("N132").Value = -345827.33

